$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.933.61'
$ws.Range("E2").Value = '  +0.90%  '

$ws.Range("D3").Value = '1.643.24'
$ws.Range("E3").Value = '  +1.33%  '

$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '''215.95'
$ws.Range("E5").Value = '  +0.55%  '

$ws.Range("D6").Value = '''0.5083'
$ws.Range("E6").Value = '  +0.20%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '''0.2606'
$ws.Range("E8").Value = '  +1.82%  '

$ws.Range("D9").Value = '''0.06485'
$ws.Range("E9").Value = '  +2.05%  '

$ws.Range("D10").Value = '''20.34'
$ws.Range("E10").Value = '  +5.72%  '

$ws.Range("D11").Value = '''0.07816'
$ws.Range("E11").Value = '  +0.66%  '

$ws.Range("D12").Value = '1.659.59'
$ws.Range("E12").Value = '  +2.09%  '

$ws.Range("D13").Value = '''4.276'
$ws.Range("E13").Value = '  +1.06%  '

$ws.Range("D14").Value = '1.868.46'
$ws.Range("E14").Value = '  +1.25%  '

$ws.Range("D15").Value = '''0.5687'
$ws.Range("E15").Value = '  +2.72%  '

$ws.Range("D16").Value = '0.0₅7718'
$ws.Range("E16").Value = '  +2.80%  '

$ws.Range("D17").Value = '''63.71'
$ws.Range("E17").Value = '  +0.39%  '

$ws.Range("D18").Value = '25.938.35'
$ws.Range("E18").Value = '  +0.83%  '

$ws.Range("E19").Value = '  +0.02%  '

$ws.Range("D20").Value = '''195.27'
$ws.Range("E20").Value = '  +1.04%  '

$ws.Range("D21").Value = '''4.414'
$ws.Range("E21").Value = '  +1.11%  '

$ws.Range("D22").Value = '''10.01'
$ws.Range("E22").Value = '  +2.78%  '

$ws.Range("D23").Value = '''6.268'
$ws.Range("E23").Value = '  +5.41%  '

$ws.Range("D24").Value = '''1.005'
$ws.Range("E24").Value = '  +0.10%  '

$ws.Range("E25").Value = '  -4.84%  '

$ws.Range("D26").Value = '''138.72'
$ws.Range("E26").Value = '  -1.26%  '

$ws.Range("D27").Value = '''0.1234'
$ws.Range("E27").Value = '  -1.13%  '

$ws.Range("D28").Value = '''6.888'
$ws.Range("E28").Value = '  +2.60%  '

$ws.Range("D29").Value = '''15.63'
$ws.Range("E29").Value = '  +1.51%  '

$ws.Range("E30").Value = '  +1.00%  '

$ws.Range("D31").Value = '''0.05037'
$ws.Range("E31").Value = '  +3.84%  '

$ws.Range("D32").Value = '''3.326'
$ws.Range("E32").Value = '  +1.04%  '

$ws.Range("D33").Value = '''3.273'
$ws.Range("E33").Value = '  +3.37%  '

$ws.Range("D34").Value = '''1.586'
$ws.Range("E34").Value = '  +3.03%  '

$ws.Range("D35").Value = '''2.385'
$ws.Range("E35").Value = '  +0.79%  '

$ws.Range("D36").Value = '''0.9136'
$ws.Range("E36").Value = '  +2.62%  '

$ws.Range("D37").Value = '''2.588'
$ws.Range("E37").Value = '  +2.12%  '

$ws.Range("D38").Value = '''0.5555'
$ws.Range("E38").Value = '  +1.74%  '

$ws.Range("D39").Value = '1.132.01'
$ws.Range("E39").Value = '  +0.71%  '

$ws.Range("D40").Value = '''0.01583'
$ws.Range("E40").Value = '  +1.79%  '

$ws.Range("D41").Value = '''0.9962'
$ws.Range("E41").Value = '  -0.63%  '


$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '''5.515'
$ws.Range("E42").Value = '  -1.05%  '

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '''100.13'
$ws.Range("E43").Value = '  +3.26%  '

$ws.Range("E44").Value = '  +1.14%  '

$ws.Range("D45").Value = '0.0₈111'
$ws.Range("E45").Value = '  -0.68%  '

$ws.Range("D46").Value = '''55.91'
$ws.Range("E46").Value = '  +2.50%  '

$ws.Range("D47").Value = '''0.4234'
$ws.Range("E47").Value = '  -4.18%  '

$ws.Range("D48").Value = '''7.735'
$ws.Range("E48").Value = '  +2.60%  '

$ws.Range("D49").Value = '''0.05050'
$ws.Range("E49").Value = '  -1.21%  '

$ws.Range("D50").Value = '''1.005'
$ws.Range("E50").Value = '  +0.89%  '

$ws.Range("E51").Value = '  +0.02%  '
